$wb = $excel.ActiveWorkbook

# --- Sheet "parsed mile posts" ---
$ws1 = $wb.Worksheets.Item("parsed mile posts")

# Update data values
$ws1.Range("C2").Value = 11.56
$ws1.Range("D2").Value = 50000
$ws1.Range("G2").Value = 1

# Update view: zoom and selection
$ws1.Activate()
$excel.ActiveWindow.Zoom = 90
$ws1.Range("H8").Select()

# --- Sheet "definitions" ---
$ws2 = $wb.Worksheets.Item("definitions")
$ws2.Activate()
$excel.ActiveWindow.Zoom = 90
$ws2.Range("B12").Select()

# Reactivate first sheet since tabSelected="true" is on sheet1
$ws1.Activate()
